# Auto-generated script to apply cryptos.xlsx cell-value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.784.96"
$ws.Range("E2").Value = "  +4.48%  "
$ws.Range("D3").Value = "2.775.75"
$ws.Range("E3").Value = "  +5.30%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.129"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.11%  "
$ws.Range("D15").Value = "3.215.10"
$ws.Range("E15").Value = "  +5.45%  "
$ws.Range("D16").Value = "2.774.04"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "51.759.08"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.25%  "
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  +6.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0823"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0360"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.77%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.92%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.114"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").Value = "2.086.35"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("E51").Value = "  -0.60%  "
